
# Updated cryptos list on Fri May 26 23:22:40 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto ranking table on Sheet1, rows 2-51.
#
# Some Price values look like plain decimals (e.g. "309.56") which Excel
# would otherwise auto-convert to a number on assignment - silently
# dropping significant trailing zeros (e.g. "2.010" -> 2.01) and changing
# the cell's stored type from text to numeric. To keep those cells as text
# (matching the source data, which sometimes has two separators, e.g.
# "26.898.15"), those values are entered with a leading apostrophe
# (forcing literal text, the same trick used when typing in the Excel UI)
# and the style is then reset to "Normal" so no stray numeric/text format
# is left applied to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.898.15'
$ws.Range('E2').Value = '  +1.30%  '
$ws.Range('D3').Value = '1.841.68'
$ws.Range('E3').Value = '  +1.55%  '
$ws.Range('E4').Value = '  +0.50%  '
$ws.Range('D5').Value = "'309.56"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.24%  '
$ws.Range('E6').Value = '  +0.47%  '
$ws.Range('D7').Value = "'0.4703"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.37%  '
$ws.Range('E8').Value = '  +1.57%  '
$ws.Range('D9').Value = "'0.07151"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.46%  '
$ws.Range('D10').Value = "'0.9195"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.01%  '
$ws.Range('E11').Value = '  +0.92%  '
$ws.Range('D12').Value = "'0.07628"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.02%  '
$ws.Range('D13').Value = '1.837.57'
$ws.Range('E13').Value = '  +1.17%  '
$ws.Range('D14').Value = "'5.283"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.50%  '
$ws.Range('D15').Value = "'6.403"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.69%  '
$ws.Range('D16').Value = "'88.06"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.96%  '
$ws.Range('E17').Value = '  +0.48%  '
$ws.Range('D18').Value = "'0.000008627"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.91%  '
$ws.Range('E19').Value = '  +0.40%  '
$ws.Range('D20').Value = '26.933.22'
$ws.Range('E20').Value = '  +1.33%  '
$ws.Range('D21').Value = "'14.48"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.29%  '
$ws.Range('D22').Value = "'5.016"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.16%  '
$ws.Range('E23').Value = '  +0.71%  '
$ws.Range('D24').Value = "'1.927"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.31%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('E26').Value = '  +2.06%  '
$ws.Range('D27').Value = "'2.010"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.72%  '
$ws.Range('D28').Value = "'114.26"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.70%  '
$ws.Range('D29').Value = "'4.858"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.60%  '
$ws.Range('D30').Value = "'0.08820"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.17%  '
$ws.Range('D31').Value = "'3.221"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.50%  '
$ws.Range('D32').Value = "'1.173"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.66%  '
$ws.Range('D33').Value = "'0.7447"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.54%  '
$ws.Range('D34').Value = "'2.751"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.34%  '
$ws.Range('D35').Value = "'4.474"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.11%  '
$ws.Range('D36').Value = "'1.089"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.79%  '
$ws.Range('D37').Value = "'0.01943"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.34%  '
$ws.Range('D38').Value = "'0.05233"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.06%  '
$ws.Range('D39').Value = "'2.967"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.88%  '
$ws.Range('D40').Value = "'0.5190"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.74%  '
$ws.Range('D41').Value = "'6.961"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.65%  '
$ws.Range('D42').Value = "'0.1512"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.18%  '
$ws.Range('D43').Value = "'8.158"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.62%  '
$ws.Range('D44').Value = "'10.49"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.57%  '
$ws.Range('D45').Value = "'0.4702"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.27%  '
$ws.Range('D46').Value = "'1.008"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.56%  '
$ws.Range('D47').Value = "'101.97"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.84%  '
$ws.Range('D48').Value = "'1.595"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.88%  '
$ws.Range('D49').Value = "'65.06"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.04%  '
$ws.Range('D50').Value = "'0.06032"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.66%  '
$ws.Range('D51').Value = "'0.8860"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.53%  '
